$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") values for rows 2-27 change from 45185 (2023-09-16)
# to 45204 (2023-10-05). Update each cell's underlying numeric (date) value.
for ($row = 2; $row -le 27; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
